$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 ("Light forest") gains enemiesDay / enemiesNight entries
$ws.Range("B11").Value = "Fox,Wolf"
$ws.Range("C11").Value = "Fox,Wolf,Wild boar,Fairy"

# Selection moved from D13 to C15
$ws.Range("C15").Select()
